# Auto-generated Excel COM-interop script
# Applies cached market-data value updates across all 8 leve-profit worksheets
# as produced by the scheduled market data refresh run.

$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 156.66667
$ws.Range("I9").Value = 147.5
$ws.Range("J9").Value = 175
$ws.Range("K9").Value = 147.5
$ws.Range("L9").Value = 175
$ws.Range("M9").Value = 21.5
$ws.Range("N9").Value = -513
$ws.Range("H43").Value = 6565
$ws.Range("I43").Value = 5495.5
$ws.Range("J43").Value = 6870.5713
$ws.Range("K43").Value = 5495.5
$ws.Range("L43").Value = 6870.5713
$ws.Range("M43").Value = -5426.5
$ws.Range("N43").Value = -7008.5713
$ws.Range("H64").Value = 10001.875
$ws.Range("I64").Value = 10007.5
$ws.Range("K64").Value = 10007.5
$ws.Range("M64").Value = -9759.5
$ws.Range("H67").Value = 10001.875
$ws.Range("I67").Value = 10007.5
$ws.Range("K67").Value = 10007.5
$ws.Range("M67").Value = -9149.5
$ws.Range("H88").Value = 1000
$ws.Range("I88").Value = 1000
$ws.Range("K88").Value = 1000
$ws.Range("M88").Value = -594
$ws.Range("H91").Value = 1000
$ws.Range("I91").Value = 1000
$ws.Range("K91").Value = 1000
$ws.Range("M91").Value = 404
$ws.Range("H93").Value = 49999.5
$ws.Range("J93").Value = 49999.5
$ws.Range("L93").Value = 49999.5
$ws.Range("N93").Value = -54991.5
$ws.Range("H106").Value = 200001000
$ws.Range("J106").Value = 1000
$ws.Range("L106").Value = 1000
$ws.Range("N106").Value = -2262
$ws.Range("H125").Value = 1005
$ws.Range("J125").Value = 1036
$ws.Range("L125").Value = 9324
$ws.Range("N125").Value = -14244
$ws.Range("H138").Value = 3029.1072
$ws.Range("I138").Value = 1019
$ws.Range("J138").Value = 3466.087
$ws.Range("K138").Value = 3057
$ws.Range("L138").Value = 10398.261
$ws.Range("M138").Value = 2083
$ws.Range("N138").Value = -20678.261

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 981.1667
$ws.Range("I35").Value = 981.1667
$ws.Range("K35").Value = 981.1667
$ws.Range("M35").Value = -575.1667
$ws.Range("H74").Value = 1991.1666
$ws.Range("I74").Value = 1991.1666
$ws.Range("K74").Value = 1991.1666
$ws.Range("M74").Value = -1117.1666
$ws.Range("H77").Value = 1991.1666
$ws.Range("I77").Value = 1991.1666
$ws.Range("K77").Value = 9955.833000000001
$ws.Range("M77").Value = -5587.833000000001
$ws.Range("H102").Value = 70000280
$ws.Range("I102").Value = 5000415
$ws.Range("K102").Value = 5000415
$ws.Range("M102").Value = -4998793
$ws.Range("H122").Value = 2949
$ws.Range("I122").Value = 2949
$ws.Range("K122").Value = 8847
$ws.Range("M122").Value = -6397

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 18995
$ws.Range("J76").Value = 18995
$ws.Range("L76").Value = 18995
$ws.Range("N76").Value = -19625
$ws.Range("H79").Value = 18995
$ws.Range("J79").Value = 18995
$ws.Range("L79").Value = 18995
$ws.Range("N79").Value = -21179
$ws.Range("H106").Value = 21924.5
$ws.Range("J106").Value = 21924.5
$ws.Range("L106").Value = 21924.5
$ws.Range("N106").Value = -24448.5
$ws.Range("H110").Value = 123333.336
$ws.Range("I110").Value = 90000
$ws.Range("J110").Value = 140000
$ws.Range("K110").Value = 90000
$ws.Range("L110").Value = 140000
$ws.Range("M110").Value = -85910
$ws.Range("N110").Value = -148180

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 248999.5
$ws.Range("J9").Value = 248999.5
$ws.Range("L9").Value = 248999.5
$ws.Range("N9").Value = -249335.5
$ws.Range("H31").Value = 6776.2
$ws.Range("I31").Value = 3999.5
$ws.Range("K31").Value = 3999.5
$ws.Range("M31").Value = -3704.5
$ws.Range("H32").Value = 4096.6
$ws.Range("I32").Value = 2161
$ws.Range("K32").Value = 2161
$ws.Range("M32").Value = -1845
$ws.Range("H33").Value = 14516.444
$ws.Range("I33").Value = 1521.2858
$ws.Range("K33").Value = 1521.2858
$ws.Range("M33").Value = -1142.2858
$ws.Range("H34").Value = 6776.2
$ws.Range("I34").Value = 3999.5
$ws.Range("K34").Value = 3999.5
$ws.Range("M34").Value = -3797.5
$ws.Range("H35").Value = 2043.5
$ws.Range("I35").Value = 2206.8572
$ws.Range("K35").Value = 2206.8572
$ws.Range("M35").Value = -1912.8572
$ws.Range("H58").Value = 1467.25
$ws.Range("I58").Value = 1462.5714
$ws.Range("K58").Value = 1462.5714
$ws.Range("M58").Value = -1259.5714
$ws.Range("H74").Value = 67125
$ws.Range("J74").Value = 67125
$ws.Range("L74").Value = 67125
$ws.Range("N74").Value = -68873
$ws.Range("H77").Value = 67125
$ws.Range("J77").Value = 67125
$ws.Range("L77").Value = 201375
$ws.Range("N77").Value = -210111
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
$ws.Range("H93").Value = 18498.572
$ws.Range("I93").Value = 11248.333
$ws.Range("K93").Value = 11248.333
$ws.Range("M93").Value = -9376.333000000001
$ws.Range("H95").Value = 30312.375
$ws.Range("J95").Value = 30312.375
$ws.Range("L95").Value = 30312.375
$ws.Range("N95").Value = -35804.375
$ws.Range("H99").Value = 1252023.6
$ws.Range("I99").Value = 1001840
$ws.Range("J99").Value = 1668996.4
$ws.Range("K99").Value = 1001840
$ws.Range("L99").Value = 1668996.4
$ws.Range("M99").Value = -1000342
$ws.Range("N99").Value = -1671992.4
$ws.Range("H103").Value = 44971.5
$ws.Range("I103").Value = 44971.5
$ws.Range("K103").Value = 44971.5
$ws.Range("M103").Value = -43799.5
$ws.Range("H122").Value = 1860
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1252023.6
$ws.Range("I126").Value = 1001840
$ws.Range("J126").Value = 1668996.4
$ws.Range("K126").Value = 3005520
$ws.Range("L126").Value = 5006989.199999999
$ws.Range("M126").Value = -3003050
$ws.Range("N126").Value = -5011929.199999999
$ws.Range("H136").Value = 1467.25
$ws.Range("I136").Value = 1462.5714
$ws.Range("K136").Value = 4387.7142
$ws.Range("M136").Value = -1837.7142
$ws.Range("H141").Value = 487505.34
$ws.Range("J141").Value = 1066665
$ws.Range("L141").Value = 1066665
$ws.Range("N141").Value = -1077025

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 362
$ws.Range("I2").Value = 753.6667
$ws.Range("K2").Value = 4522.0002
$ws.Range("M2").Value = -4409.0002
$ws.Range("H4").Value = 58570.113
$ws.Range("I4").Value = 1378.8
$ws.Range("K4").Value = 4136.4
$ws.Range("M4").Value = -4024.4

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 874.5
$ws.Range("I102").Value = 589.4
$ws.Range("K102").Value = 589.4
$ws.Range("M102").Value = 1032.6
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72080

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
$ws.Range("H136").Value = 3279
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 2798.3333
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 8394.999899999999
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = -13494.9999

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3315.8333
$ws.Range("I62").Value = 4000.6667
$ws.Range("K62").Value = 4000.6667
$ws.Range("M62").Value = -3376.6667
$ws.Range("H63").Value = 28428.572
$ws.Range("J63").Value = 28428.572
$ws.Range("L63").Value = 28428.572
$ws.Range("N63").Value = -29676.572
$ws.Range("H65").Value = 3315.8333
$ws.Range("I65").Value = 4000.6667
$ws.Range("K65").Value = 20003.3335
$ws.Range("M65").Value = -16883.3335
$ws.Range("H66").Value = 28428.572
$ws.Range("J66").Value = 28428.572
$ws.Range("L66").Value = 85285.716
$ws.Range("N66").Value = -91525.716
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

